# Update the predicted values (column B) in Sheet1 with the new,
# non-depurated (undebugged) predictions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 31760.74690414173
$ws.Range("B3").Value = 32258.02771795394
$ws.Range("B4").Value = 32772.06339788926
$ws.Range("B5").Value = 33334.29282539665
$ws.Range("B6").Value = 33522.81934132759
$ws.Range("B7").Value = 33773.75764673101
